$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.877
$ws.Range("B10").Value = 6.348
$ws.Range("B12").Value = 5.553
$ws.Range("D12").Value = -6.985000000000001
$ws.Range("D17").Value = -8.282
$ws.Range("B18").Value = 5.194
$ws.Range("D26").Value = -7.195
$ws.Range("D27").Value = -7.935
$ws.Range("D28").Value = -8.047999999999998
$ws.Range("B37").Value = 8.882000000000001
$ws.Range("D37").Value = -7.900999999999999
$ws.Range("B55").Value = 4.763
$ws.Range("D65").Value = -7.67
$ws.Range("B68").Value = 5.084999999999999
$ws.Range("D73").Value = -8.129999999999999
$ws.Range("B77").Value = 5.766
$ws.Range("B78").Value = 7.507000000000001
$ws.Range("D84").Value = -8.4
$ws.Range("D85").Value = -8.722
$ws.Range("D93").Value = -7.007000000000001
$ws.Range("D95").Value = -7.569
$ws.Range("D98").Value = -7.231
$ws.Range("D99").Value = -8.189
$ws.Range("D101").Value = -8.040000000000001
